# Generate Report for handback
# Updates the localization-status workbook to reflect a completed handback:
#  - Status changes from "Ready for handoff" to "Handed back: in sync with en-US"
#  - Latest Target File / Latest Handback File columns (E/F) are populated with
#    the same file references as the handoff columns (A/C), each as hyperlinks
#  - Latest Handback DateTime (G) is stamped with the handback timestamp

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# 1. Update the Status column everywhere "Ready for handoff" currently shows.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Cells.Item(2, 2).Value2 = $newStatus
$wsOverview.Cells.Item(2, 3).Value2 = $newStatus
$wsOverview.Cells.Item(3, 2).Value2 = $newStatus
$wsOverview.Cells.Item(3, 3).Value2 = $newStatus

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Cells.Item(2, 2).Value2 = $newStatus
$wsZh.Cells.Item(3, 2).Value2 = $newStatus

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Cells.Item(2, 2).Value2 = $newStatus
$wsDe.Cells.Item(3, 2).Value2 = $newStatus

# ---------------------------------------------------------------------------
# Helper: copy the blue-underline "hyperlink" look of column A / C cells.
# ---------------------------------------------------------------------------
function Set-HandbackLink($ws, $row, $col, $url, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $ws.Hyperlinks.Add($cell, $url, [System.Type]::Missing, [System.Type]::Missing, $text) | Out-Null
    $cell.Font.Underline = 2
    $cell.Font.Color = 15570276
}

# ---------------------------------------------------------------------------
# 2. zh-cn sheet: populate Latest Target File (E) / Latest Handback File (F)
#    with the same files as the handoff, and stamp the handback datetime (G).
# ---------------------------------------------------------------------------
Set-HandbackLink $wsZh 2 5 "https://github.com/OpenLocalizationTest/oltest/blob/813e24a9bd2d64d7a165d7db9260070df0645306/e2e/2ac41cb1-a240-442c-a524-4efa10ef2ca6.md" "2ac41cb1-a240-442c-a524-4efa10ef2ca6.md"
Set-HandbackLink $wsZh 2 6 "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a0a0720ce4959c988dafb13560021c222c14a2f2/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/qimu/2ac41cb1-a240-442c-a524-4efa10ef2ca6.df4bfcc7f35b9f2688b33325409ea6aadb12d18f.zh-cn.xlf" "2ac41cb1-a240-442c-a524-4efa10ef2ca6.df4bfcc7f35b9f2688b33325409ea6aadb12d18f.zh-cn.xlf"
$wsZh.Cells.Item(2, 7).Value2 = "2016-01-25 03:37:17"

Set-HandbackLink $wsZh 3 5 "https://github.com/OpenLocalizationTest/oltest/blob/813e24a9bd2d64d7a165d7db9260070df0645306/e2e/340eefdd-c01c-4f44-96d8-19a1448a7aab.md" "340eefdd-c01c-4f44-96d8-19a1448a7aab.md"
Set-HandbackLink $wsZh 3 6 "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a0a0720ce4959c988dafb13560021c222c14a2f2/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/qimu/340eefdd-c01c-4f44-96d8-19a1448a7aab.05a2c2c3a444464970b708d766c35f4ba88bb0ce.zh-cn.xlf" "340eefdd-c01c-4f44-96d8-19a1448a7aab.05a2c2c3a444464970b708d766c35f4ba88bb0ce.zh-cn.xlf"
$wsZh.Cells.Item(3, 7).Value2 = "2016-01-25 03:37:17"

# ---------------------------------------------------------------------------
# 3. de-de sheet: same idea, with the de-de handoff files and timestamp.
# ---------------------------------------------------------------------------
Set-HandbackLink $wsDe 2 5 "https://github.com/OpenLocalizationTest/oltest/blob/813e24a9bd2d64d7a165d7db9260070df0645306/e2e/2ac41cb1-a240-442c-a524-4efa10ef2ca6.md" "2ac41cb1-a240-442c-a524-4efa10ef2ca6.md"
Set-HandbackLink $wsDe 2 6 "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3b6dc93023672594488cbc93afbe4ce0d2624122/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/qimu/2ac41cb1-a240-442c-a524-4efa10ef2ca6.df4bfcc7f35b9f2688b33325409ea6aadb12d18f.de-de.xlf" "2ac41cb1-a240-442c-a524-4efa10ef2ca6.df4bfcc7f35b9f2688b33325409ea6aadb12d18f.de-de.xlf"
$wsDe.Cells.Item(2, 7).Value2 = "2016-01-25 03:37:33"

Set-HandbackLink $wsDe 3 5 "https://github.com/OpenLocalizationTest/oltest/blob/813e24a9bd2d64d7a165d7db9260070df0645306/e2e/340eefdd-c01c-4f44-96d8-19a1448a7aab.md" "340eefdd-c01c-4f44-96d8-19a1448a7aab.md"
Set-HandbackLink $wsDe 3 6 "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3b6dc93023672594488cbc93afbe4ce0d2624122/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/qimu/340eefdd-c01c-4f44-96d8-19a1448a7aab.05a2c2c3a444464970b708d766c35f4ba88bb0ce.de-de.xlf" "340eefdd-c01c-4f44-96d8-19a1448a7aab.05a2c2c3a444464970b708d766c35f4ba88bb0ce.de-de.xlf"
$wsDe.Cells.Item(3, 7).Value2 = "2016-01-25 03:37:33"

Write-Host "Handback report generated"
